$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("razaoSocial"), shifting B:P to C:Q
$ws.Range("B1:B2").EntireColumn.Insert()

# New column B inherits formatting oddly from the insert; re-sync formats
# from the neighboring (already shifted) cells so the look matches the rest
# of the header row / data row.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column B with header "id" and value "1306"
# (kept textual, matching the source export's inline-string cell)
$ws.Range("B1").Value = "id"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1306"
$ws.Range("B2").Style = "Normal"

# Update changed values after the shift
$ws.Range("L2").Value = "SN"

# idClienteIXC (now column O) must stay textual ("117695"), not numeric.
$ws.Range("O2").NumberFormat = "@"
$ws.Range("O2").Value = "117695"
$ws.Range("O2").Style = "Normal"

$ws.Range("Q2").Value = "{'type': 'success', 'message': 'Registro inserido com sucesso!', 'id': '117695', 'atualiza_campos': [{'tipo': 'r', 'campo': 'ativo', 'valor': 'S'}, {'tipo': 'i', 'campo': 'data_cadastro', 'valor': '19/03/2025'}, {'tipo': 'i', 'campo': 'filial_id', 'valor': '35'}, {'tipo': 'i', 'campo': 'latitude', 'valor': ''}, {'tipo': 'i', 'campo': 'longitude', 'valor': ''}, {'tipo': 'i', 'campo': 'id_conta', 'valor': '919706'}, {'tipo': 'd', 'campo': 'crm_data_vencemos', 'valor': ''}, {'tipo': 'r', 'campo': 'convert_cliente_forn', 'valor': ''}, {'tipo': 'd', 'campo': 'crm_data_perdemos', 'valor': ''}, {'tipo': 'i', 'campo': 'crm_data_sem_viabilidade', 'valor': ''}, {'tipo': 'i', 'campo': 'crm_data_sem_porta_disponivel', 'valor': ''}, {'tipo': 'i', 'campo': 'crm_data_abortamos', 'valor': ''}, {'tipo': 'i', 'campo': 'crm_data_negociando', 'valor': ''}, {'tipo': 'i', 'campo': 'crm_data_apresentando', 'valor': ''}, {'tipo': 'i', 'campo': 'crm_data_sondagem', 'valor': ''}, {'tipo': 'i', 'campo': 'crm_data_novo', 'valor': ''}]}"
